$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.371.59"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.868.11"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7033"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07918"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07837"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.59%  "
$ws.Range("D12").Value = "1.859.89"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.178"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008387"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "29.364.75"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.78%  "
$ws.Range("D20").Value = "2.114.20"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.646"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1552"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.001"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.313"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.246"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05275"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.893"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.174"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7479"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "1.272.11"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.768"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8924"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.87%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  -4.83%  "
$ws.Range("D47").Value = "2.013.69"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.619"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.798"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5181"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4298"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
